# Apply crypto price/volume updates to match target commit (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '37.368.34'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '2.030.04'
$ws.Range('E3').Value = '  +3.14%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.23'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.51'
$ws.Range('E7').Value = '  -2.36%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.390'
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('E10').Value = '  +2.12%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.00'
$ws.Range('E12').Value = '  +5.41%  '
$ws.Range('D13').Value = '2.333.10'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.839'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.73'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').Value = '2.034.21'
$ws.Range('E17').Value = '  +3.05%  '
$ws.Range('D18').Value = '37.291.55'
$ws.Range('E18').Value = '  +2.08%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.20'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').Value = '0.0₃0858'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '228.57'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.55'
$ws.Range('E24').Value = '  +4.61%  '
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.25'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '163.56'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('E28').Value = '  -5.29%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.87'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.78'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0670'
$ws.Range('E33').Value = '  +9.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.56'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  +9.00%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.56'
$ws.Range('E36').Value = '  +6.14%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.45'
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('E40').Value = '  +3.42%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0972'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('E42').Value = '  +4.28%  '
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  +4.41%  '
$ws.Range('D45').Value = '1.398.48'
$ws.Range('E45').Value = '  +2.33%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.29'
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('E47').Value = '  +2.58%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.47'
$ws.Range('E48').Value = '  +3.91%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.11'
$ws.Range('E49').Value = '  +14.33%  '
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('D51').Value = '2.224.21'
$ws.Range('E51').Value = '  +3.25%  '
